$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ipaddress row (row 2) with the new IP address
$ws.Range("B2").Value = "192.168.122.1"

# Update the hostname row (row 7) with the new hostname
$ws.Range("B7").Value = "adeye06u"

# Move the selection to B2, matching the saved cursor position
$ws.Range("B2").Select()
